$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 36
$ws.Cells.Item($row, 1).Value = "Record"
$ws.Cells.Item($row, 2).Value = "RJ Record"
$ws.Cells.Item($row, 3).Value = "Trânsito"
$ws.Cells.Item($row, 4).Value = "2025-04-01T18:17"
$ws.Cells.Item($row, 5).Value = "Neutro"
$ws.Cells.Item($row, 6).Value = "Traído pelo GPS. Carreta atinge poste e deixa Centro de Campos sem energia. Motorista teria errado o caminho e acabou subindo no calçadão do Centro. Imagens no Boulevard Francisco de Paula Carneiro, no Centro. Sem energia. Entrevista com trabalhadores da área central. Motorista vinha de Santa Catarina com destino a São Pedro da Aldeia. Entrevista com motorista da carreta. Área foi isolada. Entrevista com comerciante. *matéria* Repórter *ao vivo* do local com atualizações. Também foi veiculada no Balanço Geral. "
